$d = $word.ActiveDocument
$range = $d.Content
$found1 = $range.Find.Execute("gitignore in use.", $false, $false, $false, $false, $false, $true, 1, $false, "gitignore in use.", 2)
Write-Output ("Found1: " + $found1)

$range2 = $d.Content
$found2 = $range2.Find.Execute("appropriate. gitignore", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spacePos = $range2.End - 9 - 1
$spaceRange = $d.Range($spacePos, $spacePos+1)
Write-Output ("spaceRange text=[" + $spaceRange.Text + "]")
$spaceRange.Text = "XX"
$fullRange = $d.Content
Write-Output ("Text now: [" + $fullRange.Text.Substring(1280, 150) + "]")
